$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: rebuild hyperlinks for F2:F9 with the new (shifted) target URLs.
# The engine's Hyperlinks.Delete only works reliably when invoked on a Range
# that currently owns a hyperlink; it clears hyperlinks sheet-wide as a side
# effect, so we delete once and then re-add all eight in row order.
$ws.Range("F2").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5455098') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5445159') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5445154') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5459299') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5459200') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5459128') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5458992') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5459232') | Out-Null

# Step 2: write the cell values for the new row 2 and the rows shifted down from it.
$ws.Range("A2").Value = '2025-12-23 12:39:26'
$ws.Range("B2").Value = '大手SIer等のAIソリューション開発・導入を支援してくださるエンジニア・PM募集'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5455098'
$ws.Range("G2").Value = 375
$ws.Range("H2").Value = '🔥AI,Ai ◆開発'

$ws.Range("A3").Value = '2025-12-23 12:39:26'
$ws.Range("B3").Value = '法人向け生成AIサービス(RAG・議事録機能)の設計・開発を支援エンジニア募集(AI/バックエンド)'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5445159'
$ws.Range("G3").Value = 368
$ws.Range("H3").Value = '🔥AI,Ai ◆開発'

$ws.Range("A4").Value = '2025-12-23 12:39:26'
$ws.Range("B4").Value = 'B2B向け生成AIサービス(チャット・RAG)の新規開発プロジェクト推進を支援してくださるPM募集'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5445154'
$ws.Range("G4").Value = 368
$ws.Range("H4").Value = '🔥AI,Ai ◆開発'

$ws.Range("A5").Value = '2025-12-23 12:39:26'
$ws.Range("B5").Value = 'next.js環境下でstripe実装(オーソリ処理含む)の実装相談'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5459299'
$ws.Range("G5").Value = 128
$ws.Range("H5").Value = '🔥Next.js'

$ws.Range("A6").Value = '2025-12-23 12:39:26'
$ws.Range("B6").Value = '【急募】SNS技術検証アプリ開発'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5459200'
$ws.Range("G6").Value = 93
$ws.Range("H6").Value = '◆開発 ◇アプリ'

$ws.Range("A7").Value = '2025-12-23 12:39:26'
$ws.Range("B7").Value = '【急募】n8nを使った請求書自動化プロジェクトの依頼'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5459128'
$ws.Range("G7").Value = 88
$ws.Range("H7").Value = '◆自動化'

$ws.Range("A8").Value = '2025-12-23 12:39:26'
$ws.Range("B8").Value = '【急募】女性顧客向けチャットボット開発のプロを探しています!'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5458992'
$ws.Range("G8").Value = 68
$ws.Range("H8").Value = '◆開発'

$ws.Range("A9").Value = '2025-12-23 12:39:26'
$ws.Range("B9").Value = '【電卓設計】ハードウェアとソフトウェアの専門家を募集!'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5459232'
$ws.Range("G9").Value = 18
$ws.Range("H9").Value = ""

